$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.26%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.48%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.770"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.44%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08313"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.86%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.782"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.34%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.501"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.59%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.00%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.91%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9358"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.92%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1248"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.56%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1940"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.71%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09498"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.89%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03978"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "7.11%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1064"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.79%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001301"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.54%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005943"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.36%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.521"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.43%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.061"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "9.73%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1370"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.66%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2570"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04386"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.83%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.03%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004338"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.03%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.76%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003993"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.00%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.64%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05597"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.10%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007926"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.61%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1424"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.53%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009072"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.85%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002156"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.98%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009925"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-22.77%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007211"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.92%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.004032"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "14.40%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.13%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"
